$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the MYSQL connection details (columns E:G) - match the
# existing header style (centered, style index 1)
$ws.Range("E1").Value = "HOSTNAME"
$ws.Range("F1").Value = "PORTNO"
$ws.Range("G1").Value = "USERNAME"
$ws.Range("E1:G1").HorizontalAlignment = -4108
$ws.Range("E1:G1").VerticalAlignment = -4108

# Row 2 (the MYSQL row) now carries a real connection: server/db/host/port/user
$ws.Range("C2").Value = "TEST101"
$ws.Range("D2").Value = "mydb1"
$ws.Range("E2").Value = "127.0.0.1"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("F2").Value = 3306
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("G2").Value = "root"
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("G2").VerticalAlignment = -4108

# Move the view / selection the way the author left it
$ws.Range("G2").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 5
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
